{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the same textual changes described by the XML diff:\n//  1. \"the selected options are stored\" -> \"the selected options/questions are stored\"\n//     and \"written to a file and transferred off.\" -> \"written to a file, and transferred off.\"\n//  2. \"whats\" -> \"what is\"\n//  3. \"tried this I always\" -> \"tried this as I always\"\n//     and append a new sentence about reformatting the SD card / running with python3.\n\n// --- Change 1a: \"options\" -> \"options/questions\" -------------------------\nlet r1 = context.document.body.search(\"the selected options are stored\", { matchCase: true });\nr1.load(\"text\");\nawait context.sync();\n\nif (r1.items.length > 0) {\n  r1.items[0].insertText(\n    \"the selected options/questions are stored\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// --- Change 1b: \"a file and transferred\" -> \"a file, and transferred\" ---\nlet r2 = context.document.body.search(\"written to a file and transferred off.\", { matchCase: true });\nr2.load(\"text\");\nawait context.sync();\n\nif (r2.items.length > 0) {\n  r2.items[0].insertText(\n    \"written to a file, and transferred off.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// --- Change 2: \"whats\" -> \"what is\" ---------------------------------------\nlet r3 = context.document.body.search(\"whats\", { matchCase: true, matchWholeWord: true });\nr3.load(\"text\");\nawait context.sync();\n\nif (r3.items.length > 0) {\n  r3.items[0].insertText(\"what is\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 3a: \"tried this I always\" -> \"tried this as I always\" --------\nlet r4 = context.document.body.search(\"tried this I always\", { matchCase: true });\nr4.load(\"text\");\nawait context.sync();\n\nif (r4.items.length > 0) {\n  r4.items[0].insertText(\"tried this as I always\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 3b: append the new sentence after the shutdown-button text ---\nlet r5 = context.document.body.search(\n  \"and had the shutdown button modified to close the game.\",\n  { matchCase: true }\n);\nr5.load(\"text\");\nawait context.sync();\n\nif (r5.items.length > 0) {\n  r5.items[0].insertText(\n    \" When we got the project it was easiest to just reformat the whole SD card and setup the game/OS again. It is simple just make sure when executing the game from cmd you run it with the python3 command.\",\n    Word.InsertLocation.end\n  );\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the same textual changes described by the XML diff:\n#  1. \"the selected options are stored\" -> \"the selected options/questions are stored\"\n#     and \"written to a file and transferred off.\" -> \"written to a file, and transferred off.\"\n#  2. \"whats\" -> \"what is\"\n#  3. \"tried this I always\" -> \"tried this as I always\"\n#     and append a new sentence about reformatting the SD card / running with python3.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $findText\n  $find.Replacement.Text = $replaceText\n  $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# --- Change 1a: \"options\" -> \"options/questions\" --------------------------\nReplace-Text \"the selected options are stored\" \"the selected options/questions are stored\"\n\n# --- Change 1b: \"a file and transferred\" -> \"a file, and transferred\" -----\nReplace-Text \"written to a file and transferred off.\" \"written to a file, and transferred off.\"\n\n# --- Change 2: \"whats\" -> \"what is\" ----------------------------------------\nReplace-Text \"whats\" \"what is\"\n\n# --- Change 3a: \"tried this I always\" -> \"tried this as I always\" --------\nReplace-Text \"tried this I always\" \"tried this as I always\"\n\n# --- Change 3b: append the new sentence after the shutdown-button text ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"and had the shutdown button modified to close the game.\"\n$found = $find2.Execute()\nif ($found) {\n  $rng = $find2.Parent\n  $rng.Collapse(0)\n  $rng.InsertAfter(\" When we got the project it was easiest to just reformat the whole SD card and setup the game/OS again. It is simple just make sure when executing the game from cmd you run it with the python3 command.\")\n}\n"}
